$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text so values like "0.5180" keep trailing zeros
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.074.29"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "2.106.66"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "345.67"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").Value = "0.5180"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").Value = "0.09439"
$ws.Range("E9").Value = "  +3.55%  "
$ws.Range("D10").Value = "52.63"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").Value = "1.177"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "25.32"
$ws.Range("E12").Value = "  +3.45%  "
$ws.Range("D13").Value = "2.118.89"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "6.742"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "8.137"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "99.86"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "20.74"
$ws.Range("E19").Value = "  +6.08%  "
$ws.Range("D20").Value = "0.06707"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").Value = "6.235"
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("D23").Value = "30.169.91"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").Value = "12.69"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "2.336"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").Value = "2.352.32"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "22.09"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").Value = "164.40"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").Value = "2.552"
$ws.Range("D30").Value = "133.84"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").Value = "1.169"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("D32").Value = "0.1061"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").Value = "1.637"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "6.261"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "3.963"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "6.230"
$ws.Range("E36").Value = "  +4.42%  "
$ws.Range("D37").Value = "10.17"
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("D38").Value = "0.02569"
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").Value = "0.06798"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").Value = "0.2289"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "0.6988"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").Value = "12.56"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "1.313"
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("D44").Value = "0.6725"
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").Value = "14.30"
$ws.Range("E45").Value = "  -6.15%  "
$ws.Range("D46").Value = "2.289"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "3.640"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").Value = "0.00000000357"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").Value = "1.223"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").Value = "82.79"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").Value = "0.07215"
$ws.Range("E51").Value = "  -1.27%  "
